# studentsList.xlsx: ExcelDataService now only needs the student's name,
# national ID and application status, so collapse the "college / address /
# grade / percentage / date" columns (D:H) out of the sheet and turn the
# former "College" column (C) into the new "Notes/Status" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for column C switches from "College" to "Notes"
$ws.Range("C1").Value = "ملاحظات"

# New status per student (write row 3 before row 2 so the shared-string
# table is populated in the same order the real edit produced)
$ws.Range("C3").Value = "مرفوض"
$ws.Range("C2").Value = "مقبول لم يستدل "
$ws.Range("C4:C10").Value = "لم يستدل عليه"

# The old college/address/grade/percentage/date/notes data (D:H) is gone
$ws.Range("D1:H10").ClearContents()

# Drop the now-unused trailing column (H) entirely
$ws.Columns.Item(8).Delete()

# Selection moves to the new status column
$ws.Range("C2").Select() | Out-Null
